# Misc fixes; still fixing integer truncation import issue
#
# - Row 2-5 numeric corrections on "3. SPS - Generic Parameters"
#   (Hatch Distance / Hatch Angle / Layer Angle Increment / # Inner & Outer
#   Contours / Spot Compensation / Volume Offset Hatch columns).
# - Leftover cursor/selection state on each sheet from the editing session,
#   ending with "3. SPS - Generic Parameters" as the active sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "1. General": leave the cursor parked lower on the sheet ---
$wsGeneral = $wb.Worksheets.Item("1. General")
$wsGeneral.Activate()
$wsGeneral.Range("B24").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1

# --- Sheet "2. SPS - Area Specification": cursor moved to D3 ---
$wsArea = $wb.Worksheets.Item("2. SPS - Area Specification")
$wsArea.Activate()
$wsArea.Range("D3").Select()

# --- Sheet "3. SPS - Generic Parameters": fix the generic-parameter table ---
$wsGeneric = $wb.Worksheets.Item("3. SPS - Generic Parameters")
$wsGeneric.Activate()

$wsGeneric.Range("B2").Value = 0.1
$wsGeneric.Range("C2").Value = 0
$wsGeneric.Range("D2").Value = 66.7
$wsGeneric.Range("F2").Value = 2
$wsGeneric.Range("G2").Value = 2
$wsGeneric.Range("H2").Value = 1
$wsGeneric.Range("I2").Value = 0

$wsGeneric.Range("B3").Value = 0.1
$wsGeneric.Range("C3").Value = 45
$wsGeneric.Range("D3").Value = 66.7
$wsGeneric.Range("F3").Value = 2
$wsGeneric.Range("G3").Value = 2
$wsGeneric.Range("H3").Value = 1
$wsGeneric.Range("I3").Value = 0

$wsGeneric.Range("B4").Value = 0.1
$wsGeneric.Range("C4").Value = 85
$wsGeneric.Range("D4").Value = 66.7
$wsGeneric.Range("F4").Value = 2
$wsGeneric.Range("G4").Value = 2
$wsGeneric.Range("H4").Value = 1
$wsGeneric.Range("I4").Value = 0

$wsGeneric.Range("B5").Value = 0.1
$wsGeneric.Range("C5").Value = 125
$wsGeneric.Range("D5").Value = 66.7
$wsGeneric.Range("F5").Value = 2
$wsGeneric.Range("G5").Value = 2
$wsGeneric.Range("H5").Value = 1
$wsGeneric.Range("I5").Value = 0

# Final cursor position - this sheet ends up the active tab.
$wsGeneric.Range("I6").Select()
